$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("model_10_0_0", 0.419820259694951, 0.1508400527664997, 0.06568731025436791, 0.3255605975328584, 0.6420879364013672, 1.730097413063049, 0.2272701263427734, 1.022884607315063)
    ,@("model_10_0_24", 0.4969228240677905, -0.2396846936711678, -2.958831824791341, -0.1804578123918912, 0.5567582249641418, 2.525761365890503, 0.9629797339439392, 1.790334582328796)
    ,@("model_10_0_23", 0.4992061052142013, -0.2335032557486374, -2.943191963401899, -0.1748807793798524, 0.5542313456535339, 2.513166904449463, 0.9591754078865051, 1.7818763256073)
    ,@("model_10_0_22", 0.4999152102463791, -0.2335783442883954, -2.916870107171301, -0.172948193177267, 0.5534464716911316, 2.513319969177246, 0.9527726173400879, 1.778945207595825)
    ,@("model_10_0_21", 0.5148962607036258, -0.1824374153236399, -2.921112719552709, -0.1368967211696703, 0.5368668437004089, 2.409124374389648, 0.9538046717643738, 1.724268078804016)
    ,@("model_10_0_20", 0.5164138484376217, -0.1865260207426036, -2.832788886300593, -0.1331381778284531, 0.5351873636245728, 2.417454481124878, 0.9323199987411499, 1.718567490577698)
    ,@("model_10_0_19", 0.5276554921772456, -0.1470391388140471, -2.844337142768619, -0.105927266958471, 0.522746205329895, 2.337003231048584, 0.9351291656494141, 1.677298307418823)
    ,@("model_10_0_1", 0.5288798698710702, 0.3189295441339353, -1.096691514476998, 0.357374229522472, 0.5213912129402161, 1.387628078460693, 0.5100170373916626, 0.9746347069740295)
    ,@("model_10_0_18", 0.5380671359296699, -0.1128822278538824, -2.832536775327245, -0.08074372548928932, 0.511223554611206, 2.267411231994629, 0.9322587251663208, 1.639103889465332)
    ,@("model_10_0_17", 0.562202191296568, -0.04911300738385549, -2.662741530402664, -0.02257554747442869, 0.484513133764267, 2.137486219406128, 0.8909562826156616, 1.550883412361145)
    ,@("model_10_0_16", 0.5819404989476881, 0.0225000566697432, -2.70260807979344, 0.0253464196535591, 0.4626686871051788, 1.991580128669739, 0.9006537199020386, 1.478202819824219)
    ,@("model_10_0_15", 0.5970848982712313, 0.0727302707069174, -2.681978615969077, 0.06262677344953438, 0.4459083080291748, 1.889240026473999, 0.8956356048583984, 1.421661853790283)
    ,@("model_10_0_14", 0.6027181045287496, 0.09024114193297916, -2.667080148839596, 0.0762048460567204, 0.4396739900112152, 1.853562831878662, 0.8920116424560547, 1.401068687438965)
    ,@("model_10_0_13", 0.6670572757386747, 0.3073390224649396, -2.651357946646561, 0.2317914741168796, 0.3684695065021515, 1.411242961883545, 0.8881872296333313, 1.165099143981934)
    ,@("model_10_0_12", 0.694459231713451, 0.4040199255396154, -2.679648699860401, 0.2984154893126019, 0.3381435871124268, 1.214263081550598, 0.8950689435005188, 1.064054131507874)
    ,@("model_10_0_11", 0.7068957371583571, 0.4752967080705121, -2.93923054432014, 0.3295155634252815, 0.3243800699710846, 1.069042205810547, 0.9582117795944214, 1.016886472702026)
    ,@("model_10_0_2", 0.7268433593104267, 0.5866644898789335, -1.55057947113447, 0.5135294163953268, 0.3023039400577545, 0.8421390056610107, 0.6204245090484619, 0.7378028631210327)
    ,@("model_10_0_10", 0.7352852167885662, 0.5749494530138082, -2.956782391586033, 0.399063556320255, 0.2929612696170807, 0.8660074472427368, 0.9624812006950378, 0.911406934261322)
    ,@("model_10_0_9", 0.7365824728627355, 0.5799447637082356, -2.960262573006388, 0.4023534166901206, 0.2915255725383759, 0.8558299541473389, 0.9633278250694275, 0.9064173698425293)
    ,@("model_10_0_8", 0.7456997853790047, 0.6078428121162333, -2.92992169844679, 0.4244843087619784, 0.2814354002475739, 0.7989897727966309, 0.9559473991394043, 0.8728526830673218)
    ,@("model_10_0_7", 0.747254858339292, 0.6105217449131943, -2.902218420987799, 0.4284807020312199, 0.2797144055366516, 0.7935316562652588, 0.949208676815033, 0.8667914867401123)
    ,@("model_10_0_6", 0.7665514841132238, 0.6550190178298698, -2.676953969625329, 0.4771292039649689, 0.2583587169647217, 0.7028719782829285, 0.894413411617279, 0.7930090427398682)
    ,@("model_10_0_5", 0.7816042624276021, 0.691639924214684, -2.494555289977589, 0.5169405941608965, 0.2416997104883194, 0.6282597184181213, 0.8500452041625977, 0.7326292991638184)
    ,@("model_10_0_3", 0.7867403416251588, 0.7112173110684303, -1.788497607859882, 0.5841542122836794, 0.2360156029462814, 0.5883723497390747, 0.6782978177070618, 0.6306900978088379)
    ,@("model_10_0_4", 0.7915050457831647, 0.7106854003209669, -2.205712927260333, 0.5522861803255941, 0.230742484331131, 0.5894560217857361, 0.7797847986221313, 0.6790226697921753)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    for ($c = 0; $c -lt 8; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $rowData[$c + 1]
    }
}
